{"js": "// Removed the 6th RQ (research question) and tweaked the wording of the\n// 4th RQ (\"...best suited\" -> \"...best suited,\") to read correctly once the\n// list item that followed it is removed.\n\nconst body = context.document.body;\n\n// --- 1. \"Which CNN architecture will be best suited\" -> \"...suited,\" ---\n// There are two occurrences of the word \"suited\" in the research-question\n// list; the one we need belongs to the paragraph that talks about the CNN\n// architecture, so disambiguate using the paragraph text before editing.\nconst suitedMatches = body.search(\"suited\", { matchCase: false, matchWholeWord: true });\nsuitedMatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < suitedMatches.items.length; i++) {\n  const match = suitedMatches.items[i];\n  const paraRange = match.paragraphs.getFirst();\n  paraRange.load(\"text\");\n  suitedMatches.items[i].paraRange = paraRange;\n}\nawait context.sync();\n\nlet cnnMatch = null;\nfor (let i = 0; i < suitedMatches.items.length; i++) {\n  const match = suitedMatches.items[i];\n  if (match.paraRange.text.indexOf(\"CNN architecture\") !== -1) {\n    cnnMatch = match;\n    break;\n  }\n}\n\nif (cnnMatch) {\n  cnnMatch.insertText(\",\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// --- 2. Remove the 6th research question paragraph entirely ---\nconst rqMatches = body.search(\n  \"Will morphological image processing operators such as Erosion and Dilation improve the classifiers\",\n  { matchCase: false }\n);\nrqMatches.load(\"items\");\nawait context.sync();\n\nif (rqMatches.items.length > 0) {\n  const rqParagraph = rqMatches.items[0].paragraphs.getFirst();\n  rqParagraph.delete();\n  await context.sync();\n}\n", "ps1": "# Removed the 6th RQ (research question) and tweaked the wording of the\n# 4th RQ (\"...best suited\" -> \"...best suited,\") to read correctly once the\n# list item that followed it is removed.\n\n$d = $word.ActiveDocument\n\n# --- 1. \"Which CNN architecture will be best suited\" -> \"...suited,\" ---\n# Scope the Find to the paragraph that actually talks about the CNN\n# architecture (the word \"suited\" also appears earlier in the RQ list), then\n# insert a comma right after it.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*CNN architecture*\") {\n        $pr = $p.Range\n        $found = $pr.Find.Execute(\"suited\")\n        if ($found) {\n            $pr.Collapse(0)  # wdCollapseEnd\n            $pr.InsertAfter(\",\")\n        }\n        break\n    }\n}\n\n# --- 2. Remove the 6th research question paragraph entirely ---\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Will morphological image processing operators*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
